# Updates the cryptocurrency price/volume table (columns D and E) on the active
# worksheet to reflect the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.477.90"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.677.18"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'217.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'0.5323"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.2695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "'0.06410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").Value = "'21.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.85%  "
$ws.Range("D11").Value = "'0.07804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "1.680.09"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").Value = "'4.516"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "'0.5592"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "0.0₅8328"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'65.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "26.509.74"
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'4.790"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").Value = "'193.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").Value = "'10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "'6.338"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'142.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'0.1282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("D26").Value = "'7.421"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").Value = "'0.06286"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.27%  "
$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'3.620"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("D32").Value = "'3.461"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "'1.694"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "'1.010"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "'0.6157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.38%  "
$ws.Range("D36").Value = "'2.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'2.787"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'6.168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.42%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").Value = "1.098.25"
$ws.Range("E40").Value = "  +6.21%  "
$ws.Range("D41").Value = "'0.8626"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "'0.9997"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'100.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Value = "1.821.88"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").Value = "'57.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "'8.144"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "'0.9992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "'0.05209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "'1.480"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.89%  "
$ws.Range("D51").Value = "'6.058"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
